$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates derived from the crypto price refresh.
# Each entry: Row number, Column letter, new text Value, and whether the
# value needs to be forced to remain a text string (some price figures are
# plain decimal numbers and Excel would otherwise auto-convert them to
# numeric cells, losing the literal text formatting used throughout column D).
$changes = @(
    @{Row=2; Col='D'; Value='46.180.34'; ForceText=$false},
    @{Row=2; Col='E'; Value='  +0.82%  '; ForceText=$false},
    @{Row=3; Col='D'; Value='2.594.37'; ForceText=$false},
    @{Row=3; Col='E'; Value='  +9.23%  '; ForceText=$false},
    @{Row=4; Col='D'; Value='0.999'; ForceText=$true},
    @{Row=4; Col='E'; Value='  -0.18%  '; ForceText=$false},
    @{Row=5; Col='D'; Value='308.24'; ForceText=$true},
    @{Row=5; Col='E'; Value='  +2.72%  '; ForceText=$false},
    @{Row=6; Col='D'; Value='100.19'; ForceText=$true},
    @{Row=6; Col='E'; Value='  +2.45%  '; ForceText=$false},
    @{Row=7; Col='D'; Value='0.596'; ForceText=$true},
    @{Row=7; Col='E'; Value='  +5.89%  '; ForceText=$false},
    @{Row=8; Col='E'; Value='  -0.07%  '; ForceText=$false},
    @{Row=9; Col='D'; Value='0.580'; ForceText=$true},
    @{Row=9; Col='E'; Value='  +14.40%  '; ForceText=$false},
    @{Row=10; Col='D'; Value='38.51'; ForceText=$true},
    @{Row=10; Col='E'; Value='  +12.76%  '; ForceText=$false},
    @{Row=11; Col='D'; Value='0.0839'; ForceText=$true},
    @{Row=12; Col='D'; Value='8.28'; ForceText=$true},
    @{Row=12; Col='E'; Value='  +16.92%  '; ForceText=$false},
    @{Row=13; Col='D'; Value='2.983.65'; ForceText=$false},
    @{Row=13; Col='E'; Value='  +9.05%  '; ForceText=$false},
    @{Row=14; Col='E'; Value='  +1.42%  '; ForceText=$false},
    @{Row=15; Col='D'; Value='2.591.49'; ForceText=$false},
    @{Row=15; Col='E'; Value='  +9.76%  '; ForceText=$false},
    @{Row=16; Col='D'; Value='0.902'; ForceText=$true},
    @{Row=16; Col='E'; Value='  +10.11%  '; ForceText=$false},
    @{Row=17; Col='D'; Value='14.82'; ForceText=$true},
    @{Row=17; Col='E'; Value='  +8.36%  '; ForceText=$false},
    @{Row=18; Col='D'; Value='46.296.20'; ForceText=$false},
    @{Row=18; Col='E'; Value='  +1.15%  '; ForceText=$false},
    @{Row=19; Col='E'; Value='  +6.72%  '; ForceText=$false},
    @{Row=20; Col='D'; Value='12.93'; ForceText=$true},
    @{Row=20; Col='E'; Value='  +1.69%  '; ForceText=$false},
    @{Row=21; Col='E'; Value='  +9.87%  '; ForceText=$false},
    @{Row=22; Col='D'; Value='71.14'; ForceText=$true},
    @{Row=22; Col='E'; Value='  +6.64%  '; ForceText=$false},
    @{Row=23; Col='D'; Value='253.88'; ForceText=$true},
    @{Row=23; Col='E'; Value='  +4.23%  '; ForceText=$false},
    @{Row=24; Col='D'; Value='3.03'; ForceText=$true},
    @{Row=24; Col='E'; Value='  +8.76%  '; ForceText=$false},
    @{Row=25; Col='D'; Value='2.22'; ForceText=$true},
    @{Row=25; Col='E'; Value='  +16.11%  '; ForceText=$false},
    @{Row=26; Col='E'; Value='  +33.27%  '; ForceText=$false},
    @{Row=27; Col='D'; Value='1.00'; ForceText=$true},
    @{Row=27; Col='E'; Value='  -0.05%  '; ForceText=$false},
    @{Row=28; Col='D'; Value='10.46'; ForceText=$true},
    @{Row=28; Col='E'; Value='  +7.98%  '; ForceText=$false},
    @{Row=29; Col='D'; Value='39.58'; ForceText=$true},
    @{Row=29; Col='E'; Value='  +1.86%  '; ForceText=$false},
    @{Row=30; Col='E'; Value='  +3.35%  '; ForceText=$false},
    @{Row=31; Col='B'; Value='LidoDAOToken'; ForceText=$false},
    @{Row=31; Col='C'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; ForceText=$false},
    @{Row=31; Col='D'; Value='3.73'; ForceText=$true},
    @{Row=31; Col='E'; Value='  -1.42%  '; ForceText=$false},
    @{Row=32; Col='B'; Value='Filecoin'; ForceText=$false},
    @{Row=32; Col='C'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText=$false},
    @{Row=32; Col='D'; Value='6.05'; ForceText=$true},
    @{Row=32; Col='E'; Value='  +9.57%  '; ForceText=$false},
    @{Row=33; Col='D'; Value='2.28'; ForceText=$true},
    @{Row=33; Col='E'; Value='  +18.85%  '; ForceText=$false},
    @{Row=34; Col='D'; Value='2.89'; ForceText=$true},
    @{Row=34; Col='E'; Value='  +6.13%  '; ForceText=$false},
    @{Row=35; Col='D'; Value='152.45'; ForceText=$true},
    @{Row=35; Col='E'; Value='  +3.66%  '; ForceText=$false},
    @{Row=36; Col='D'; Value='0.0831'; ForceText=$true},
    @{Row=36; Col='E'; Value='  +8.02%  '; ForceText=$false},
    @{Row=37; Col='E'; Value='  +3.29%  '; ForceText=$false},
    @{Row=38; Col='E'; Value='  +5.61%  '; ForceText=$false},
    @{Row=39; Col='D'; Value='16.01'; ForceText=$true},
    @{Row=39; Col='E'; Value='  +6.12%  '; ForceText=$false},
    @{Row=40; Col='D'; Value='4.19'; ForceText=$true},
    @{Row=40; Col='E'; Value='  +9.31%  '; ForceText=$false},
    @{Row=41; Col='D'; Value='3.61'; ForceText=$true},
    @{Row=41; Col='E'; Value='  +13.22%  '; ForceText=$false},
    @{Row=42; Col='D'; Value='0.0323'; ForceText=$true},
    @{Row=42; Col='E'; Value='  +8.36%  '; ForceText=$false},
    @{Row=43; Col='B'; Value='EnergySwap'; ForceText=$false},
    @{Row=43; Col='C'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText=$false},
    @{Row=43; Col='D'; Value='20.39'; ForceText=$true},
    @{Row=43; Col='E'; Value='  +42.63%  '; ForceText=$false},
    @{Row=44; Col='B'; Value='Maker'; ForceText=$false},
    @{Row=44; Col='C'; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; ForceText=$false},
    @{Row=44; Col='D'; Value='2.049.61'; ForceText=$false},
    @{Row=44; Col='E'; Value='  +5.33%  '; ForceText=$false},
    @{Row=45; Col='D'; Value='0.999'; ForceText=$true},
    @{Row=45; Col='E'; Value='  -0.09%  '; ForceText=$false},
    @{Row=46; Col='D'; Value='91.00'; ForceText=$true},
    @{Row=46; Col='E'; Value='  -4.08%  '; ForceText=$false},
    @{Row=47; Col='D'; Value='9.28'; ForceText=$true},
    @{Row=47; Col='E'; Value='  +9.06%  '; ForceText=$false},
    @{Row=48; Col='B'; Value='Aave'; ForceText=$false},
    @{Row=48; Col='C'; Value='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText=$false},
    @{Row=48; Col='D'; Value='109.19'; ForceText=$true},
    @{Row=48; Col='E'; Value='  +10.48%  '; ForceText=$false},
    @{Row=49; Col='B'; Value='Stacks'; ForceText=$false},
    @{Row=49; Col='C'; Value='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; ForceText=$false},
    @{Row=49; Col='D'; Value='1.77'; ForceText=$true},
    @{Row=49; Col='E'; Value='  -0.08%  '; ForceText=$false},
    @{Row=50; Col='D'; Value='0.200'; ForceText=$true},
    @{Row=50; Col='E'; Value='  +8.64%  '; ForceText=$false},
    @{Row=51; Col='D'; Value='2.843.39'; ForceText=$false},
    @{Row=51; Col='E'; Value='  +9.12%  '; ForceText=$false},
)

foreach ($item in $changes) {
    $addr = ($item.Col + $item.Row)
    $cell = $ws.Range($addr)
    if ($item.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $item.Value
}
